$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("K33").Value = 2397127
$ws.Range("H33").Value = 1826468.1
$ws.Range("I33").Value = 2397127
$ws.Range("M33").Value = -2396898

$ws.Range("H55").Value = 140.47058
$ws.Range("N55").Value = -624.85715
$ws.Range("J55").Value = 196.85715
$ws.Range("L55").Value = 196.85715
$ws.Range("M55").Value = 113
$ws.Range("I55").Value = 101
$ws.Range("K55").Value = 101

$ws.Range("L69").Value = 11465.4
$ws.Range("N69").Value = -13213.4
$ws.Range("H69").Value = 3821.8
$ws.Range("J69").Value = 3821.8

$ws.Range("L72").Value = 34396.2
$ws.Range("N72").Value = -43132.2
$ws.Range("H72").Value = 3821.8
$ws.Range("J72").Value = 3821.8

$ws.Range("J74").Value = 4057.1428
$ws.Range("H74").Value = 12016.667
$ws.Range("M74").Value = -22224
$ws.Range("L74").Value = 4057.1428
$ws.Range("K74").Value = 23160
$ws.Range("I74").Value = 23160
$ws.Range("N74").Value = -5929.1428

$ws.Range("J76").Value = 4605.8823
$ws.Range("I76").Value = 11149.625
$ws.Range("M76").Value = -10834.625
$ws.Range("K76").Value = 11149.625
$ws.Range("L76").Value = 4605.8823
$ws.Range("N76").Value = -5235.8823
$ws.Range("H76").Value = 7778.606

$ws.Range("L77").Value = 20285.714
$ws.Range("H77").Value = 12016.667
$ws.Range("N77").Value = -29645.714
$ws.Range("M77").Value = -111120
$ws.Range("K77").Value = 115800
$ws.Range("J77").Value = 4057.1428
$ws.Range("I77").Value = 23160

$ws.Range("H79").Value = 7778.606
$ws.Range("M79").Value = -10057.625
$ws.Range("I79").Value = 11149.625
$ws.Range("L79").Value = 4605.8823
$ws.Range("N79").Value = -6789.8823
$ws.Range("J79").Value = 4605.8823
$ws.Range("K79").Value = 11149.625

$ws.Range("H111").Value = 3510.7856
$ws.Range("K111").Value = 8550.332999999999
$ws.Range("L111").Value = 14100
$ws.Range("J111").Value = 4700
$ws.Range("I111").Value = 2850.111
$ws.Range("N111").Value = -20234
$ws.Range("M111").Value = -5483.332999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I2").Value = 1634.5555
$ws.Range("K2").Value = 1634.5555
$ws.Range("N2").Value = -3766.9167
$ws.Range("J2").Value = 3540.9167
$ws.Range("H2").Value = 2397.1
$ws.Range("L2").Value = 3540.9167
$ws.Range("M2").Value = -1521.5555

$ws.Range("I32").Value = 12094.811
$ws.Range("K32").Value = 12094.811
$ws.Range("N32").Value = -127741.45
$ws.Range("H32").Value = 30439.725
$ws.Range("M32").Value = -11807.811
$ws.Range("J32").Value = 127167.45
$ws.Range("L32").Value = 127167.45

$ws.Range("I63").Value = 2500
$ws.Range("J63").Value = 4954.5454
$ws.Range("K63").Value = 2500
$ws.Range("H63").Value = 4187.5
$ws.Range("N63").Value = -6326.5454
$ws.Range("L63").Value = 4954.5454
$ws.Range("M63").Value = -1814

$ws.Range("J66").Value = 4954.5454
$ws.Range("L66").Value = 24772.727
$ws.Range("M66").Value = -9068
$ws.Range("N66").Value = -31636.727
$ws.Range("H66").Value = 4187.5
$ws.Range("K66").Value = 12500
$ws.Range("I66").Value = 2500

$ws.Range("L110").Value = 1114.6
$ws.Range("J110").Value = 1114.6
$ws.Range("I110").Value = 1581.1111
$ws.Range("K110").Value = 1581.1111
$ws.Range("N110").Value = -5204.6
$ws.Range("H110").Value = 1508.2188
$ws.Range("M110").Value = 463.8888999999999

$ws.Range("N116").Value = -8128.9167
$ws.Range("K116").Value = 1634.5555
$ws.Range("J116").Value = 3540.9167
$ws.Range("M116").Value = 659.4445000000001
$ws.Range("H116").Value = 2397.1
$ws.Range("L116").Value = 3540.9167
$ws.Range("I116").Value = 1634.5555

$ws.Range("L122").Value = 8228.143199999999
$ws.Range("I122").Value = 24382.445
$ws.Range("H122").Value = 14915.0625
$ws.Range("M122").Value = -70697.33499999999
$ws.Range("J122").Value = 2742.7144
$ws.Range("N122").Value = -13128.1432
$ws.Range("K122").Value = 73147.33499999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("K3").Value = 1634.5555
$ws.Range("I3").Value = 1634.5555
$ws.Range("H3").Value = 2397.1
$ws.Range("J3").Value = 3540.9167
$ws.Range("M3").Value = -1520.5555
$ws.Range("N3").Value = -3768.9167
$ws.Range("L3").Value = 3540.9167

$ws.Range("I105").Value = 1423178.9
$ws.Range("K105").Value = 1423178.9
$ws.Range("H105").Value = 815709.7
$ws.Range("N105").Value = -9244.833500000001
$ws.Range("L105").Value = 5750.8335
$ws.Range("J105").Value = 5750.8335
$ws.Range("M105").Value = -1421431.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I31").Value = 1776.9524
$ws.Range("L31").Value = 8083.5713
$ws.Range("H31").Value = 3353.6072
$ws.Range("N31").Value = -8673.5713
$ws.Range("K31").Value = 1776.9524
$ws.Range("J31").Value = 8083.5713
$ws.Range("M31").Value = -1481.9524

$ws.Range("I34").Value = 1776.9524
$ws.Range("M34").Value = -1574.9524
$ws.Range("H34").Value = 3353.6072
$ws.Range("J34").Value = 8083.5713
$ws.Range("L34").Value = 8083.5713
$ws.Range("K34").Value = 1776.9524
$ws.Range("N34").Value = -8487.5713

$ws.Range("L56").Value = 7083.3335
$ws.Range("K56").Value = 7000
$ws.Range("N56").Value = -8773.333500000001
$ws.Range("I56").Value = 7000
$ws.Range("H56").Value = 7055.5557
$ws.Range("J56").Value = 7083.3335
$ws.Range("M56").Value = -6155

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 672.9583
$ws.Range("K44").Value = 902.25
$ws.Range("L44").Value = 2120.3865
$ws.Range("N44").Value = -2916.3865
$ws.Range("M44").Value = -504.25
$ws.Range("J44").Value = 706.7954999999999
$ws.Range("I44").Value = 300.75

$ws.Range("I48").Value = 400
$ws.Range("J48").Value = 1440
$ws.Range("L48").Value = 4320
$ws.Range("K48").Value = 1200
$ws.Range("H48").Value = 1142.8572
$ws.Range("M48").Value = -950
$ws.Range("N48").Value = -4820

$ws.Range("H49").Value = 8600
$ws.Range("J49").Value = 8600
$ws.Range("L49").Value = 25800
$ws.Range("N49").Value = -26112

$ws.Range("J113").Value = 1150
$ws.Range("M113").Value = -590
$ws.Range("H113").Value = 985.7143
$ws.Range("L113").Value = 3450
$ws.Range("K113").Value = 2760
$ws.Range("N113").Value = -7790
$ws.Range("I113").Value = 920

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 30547.85
$ws.Range("K70").Value = 39280.266
$ws.Range("N70").Value = -4890.6
$ws.Range("J70").Value = 4350.6
$ws.Range("M70").Value = -39010.266
$ws.Range("I70").Value = 39280.266
$ws.Range("L70").Value = 4350.6

$ws.Range("M73").Value = -38344.266
$ws.Range("N73").Value = -6222.6
$ws.Range("J73").Value = 4350.6
$ws.Range("I73").Value = 39280.266
$ws.Range("L73").Value = 4350.6
$ws.Range("H73").Value = 30547.85
$ws.Range("K73").Value = 39280.266

$ws.Range("H80").Value = 7145657
$ws.Range("N80").Value = -7147653
$ws.Range("J80").Value = 7145657
$ws.Range("L80").Value = 7145657
$ws.Range("K80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("M80").ClearContents()

$ws.Range("L83").Value = 35728285
$ws.Range("H83").Value = 7145657
$ws.Range("I83").Value = 0
$ws.Range("N83").Value = -35738269
$ws.Range("J83").Value = 7145657
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()

$ws.Range("K126").Value = 8520
$ws.Range("N126").Value = -28095.4121
$ws.Range("M126").Value = -6050
$ws.Range("I126").Value = 2840
$ws.Range("J126").Value = 7718.4707
$ws.Range("H126").Value = 6609.727
$ws.Range("L126").Value = 23155.4121

$ws.Range("M132").Value = -6451.1819
$ws.Range("I132").Value = 2993.7273
$ws.Range("K132").Value = 8981.1819
$ws.Range("N132").Value = -12484.4999
$ws.Range("J132").Value = 2474.8333
$ws.Range("H132").Value = 2810.5881
$ws.Range("L132").Value = 7424.499899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("L122").Value = 15049.9995
$ws.Range("I122").Value = 3615.6924
$ws.Range("H122").Value = 4288.16
$ws.Range("M122").Value = -8397.0772
$ws.Range("J122").Value = 5016.6665
$ws.Range("N122").Value = -19949.9995
$ws.Range("K122").Value = 10847.0772

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("J107").Value = 2442.7
$ws.Range("K107").Value = 3799.4208
$ws.Range("H107").Value = 1672.069
$ws.Range("N107").Value = -11168.1
$ws.Range("M107").Value = -1879.4208
$ws.Range("L107").Value = 7328.099999999999
$ws.Range("I107").Value = 1266.4736
